# Regenerate the "K" (strikeouts) column (column G) of the save-data sheet.
# The save-data pipeline now pulls actual strikeout counts ("K") instead of
# the old "Strike#" proxy, so column G is rewritten here with the
# recalculated per-appearance values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(0,0,0,1,1,0,2,1,0,2,1,2,1,0,2,1,0,1,1,1,0,3,1,1,1,1,2,1,1,1,1,4,0,1,0,0,0,0,1,1,3,1,0,2,2,0,3,2,0,2,0,0,1,1,2,0)

$firstRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
